$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / mean calculation
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -7
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = -7
$ws.Range("F8").Value = 8
$ws.Range("F10").Value = -7
